$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first five item rows (Kolibri, M1911, DGL-50, X18, ASVal).
# This shifts everything below up by 5 rows, matching the new table layout.
$ws.Range("A4:F8").EntireRow.Delete()

# Update the remaining probability values (columns C:F) to the corrected figures.
$values = @{
    4  = @(9, 10, 11, 11)
    5  = @(9, 10, 11, 11)
    6  = @(9, 10, 11, 11)
    7  = @(9, 9, 10, 10)
    8  = @(9, 9, 10, 10)
    9  = @(9, 9, 9, 9)
    10 = @(9, 9, 8, 8)
    11 = @(9, 7, 6, 6)
    12 = @(7, 8, 8, 8)
    13 = @(7, 7, 6, 6)
    14 = @(7, 6, 5, 5)
    15 = @(7, 6, 5, 5)
}

foreach ($r in $values.Keys) {
    $row = $values[$r]
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
}

# Update print setup (paper size / orientation) as reflected in the saved file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the selection to match the saved workbook state.
$ws.Range("E8").Select()
